# Replace the Victorian "key outbreaks" cluster table with the updated
# list of active outbreaks and their case counts (rows 2..41 on Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    '3153 Sacred Heart Community St Kilda Tier 1B Outbreak',
    '3175 The Bays Aged Care Facility Hastings Outbreak',
    '3600 Belvedere Aged Care Noble Park Outbreak',
    '3612 BlueCross Glengowrie Outbreak',
    '3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak',
    '3980 Arcare Keysborough Aged Care Keysborough Outbreak',
    '4075 Ferndale Gardens Aged Care Services Bayswater North Outbreak',
    '4518 Regis Aged Care Fawkner Outbreak',
    'AW Window Transport Group Depot North Geelong Outbreak',
    'Allied Pinnacle Factory Altona North Outbreak',
    'Bespoke Childcare Dingley Village Outbreak',
    'Bread Solutions Braeside Outbreak',
    'CS Square Caroline Springs Outbreak',
    'Child''s Play Early Learning Centre Tarneit Outbreak',
    'Community Kids Pascoe Vale Early Education Centre Pascoe Vale Outbreak',
    'Essential Caravans Somerton Outbreak',
    'Guardian Childcare Caulfield Outbreak',
    'Hello Fresh Warehouse Ravenhall Outbreak',
    'Inghams Enterprises Somerville Outbreak',
    'Lantmannen Unibake Australia Mordialloc Outbreak',
    'Launch Housing City Edge Crisis Accommodation South Melbourne Outbreak',
    'MacKillop Family Services Residential Facility Glenroy Outbreak',
    'Melbourne Assessment Prison West Melbourne Outbreak',
    'Melbourne Youth Justice Centre Parkville Outbreak',
    'Mildura Police Station Mildura Outbreak',
    'Monash Health Kingston Centre South 5 Outbreak',
    'Northern Health Northern Hospital Epping Emergency Department Tier 1B Outbreak',
    'Northern Health The Northern Hospital Epping Outbreak',
    'Oceania Meat Processors Laverton North Outbreak',
    'Public Housing 140 Brunswick Street Fitzroy Outbreak',
    'Shawlands Caravan Park Dandenong South Outbreak',
    'St Vincents Hospital Emergency Department Melbourne Outbreak',
    'Target Distribution Centre Truganina Outbreak',
    'The Robin Hood Inn Drouin West Outbreak',
    'The Toolshed Bar Private Event Noojee Outbreak',
    'Turosi Breakwater Outbreak',
    'Visy Recycling Springvale Outbreak',
    'Werribee Mercy Hospital Emergency Department Outbreak',
    'Western Health Footscray Hospital Emergency Department Outbreak',
    'Western Health Sunshine Hospital Emergency Department Outbreak'
)

$values = @(
    15,
    15,
    27,
    35,
    31,
    19,
    17,
    13,
    10,
    15,
    10,
    15,
    11,
    11,
    22,
    10,
    21,
    17,
    24,
    21,
    11,
    10,
    16,
    10,
    10,
    10,
    43,
    19,
    17,
    11,
    17,
    37,
    11,
    48,
    16,
    10,
    14,
    26,
    11,
    24
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
